# The presentation had slide 10 ("Problemas a Serem Resolvidos") and slide
# 11 ("Implementação do Código"). This edit removes the "Implementação do
# Código" slide entirely and moves "Problemas a Serem Resolvidos" up to
# become the 3rd slide in the deck.

$p = $ppt.ActivePresentation

# Delete slide 11 ("Implementação do Código") - it is removed from the deck.
$p.Slides.Item(11).Delete()

# Move slide 10 ("Problemas a Serem Resolvidos") to become slide 3.
$p.Slides.Item(10).MoveTo(3)
